$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: F3 -> "-" (was the microcontroladores entry, now moved out)
$ws.Range("F3").Value = "-"

# Row 4: D4 gets the microcontroladores entry, F4 becomes "-"
$ws.Range("D4").Value = "[-, -, 'MCT-3A-Microcontroladores', -]"
$ws.Range("F4").Value = "-"

# Row 6: D6 gets the microcontroladores entry, F6 becomes "-"
$ws.Range("D6").Value = "[-, -, 'MCT-3A-Microcontroladores', -]"
$ws.Range("F6").Value = "-"

# Row 7: D7 gets the microcontroladores entry, F7 becomes "-"
$ws.Range("D7").Value = "[-, -, 'MCT-3A-Microcontroladores', -]"
$ws.Range("F7").Value = "-"

# Row 8: D8 gets the microcontroladores entry (F8 was already "-")
$ws.Range("D8").Value = "[-, -, 'MCT-3A-Microcontroladores', -]"
